# Roboflow Annotation Report 7/10/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Append new data row (55), duplicating the style of the previous last row (54) ---
$ws.Range("D54:J54").Copy()
$ws.Range("D55:J56").PasteSpecial(-4122)
$ws.Rows.Item(55).RowHeight = 15.6
$ws.Rows.Item(56).RowHeight = 15.6

$ws.Range("D55").Value = 45937
$ws.Range("E55").Value = 150
$ws.Range("F55").Value = 776
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 1012
$ws.Range("J55").Value = "N/A"

# --- Resize the table/autofilter to include the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("D4:J56"))

# --- Update the view: scroll position + active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F58").Select()
